$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/volume text (e.g. "1.002", "45.02") is written
# back as literal text, matching the original inline-string cell content,
# instead of being auto-parsed into numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '28.320.61'
$ws.Range('E2').Value = '  +3.42%  '

# Row 3
$ws.Range('D3').Value = '1.820.32'
$ws.Range('E3').Value = '  +4.46%  '

# Row 4
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Value = '328.82'
$ws.Range('E5').Value = '  +2.47%  '

# Row 6
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.15%  '

# Row 7
$ws.Range('D7').Value = '0.4348'
$ws.Range('E7').Value = '  +3.68%  '

# Row 8
$ws.Range('D8').Value = '0.3687'
$ws.Range('E8').Value = '  +3.05%  '

# Row 9
$ws.Range('D9').Value = '45.02'
$ws.Range('E9').Value = '  -0.86%  '

# Row 10
$ws.Range('D10').Value = '0.07719'
$ws.Range('E10').Value = '  +4.13%  '

# Row 11
$ws.Range('D11').Value = '1.140'
$ws.Range('E11').Value = '  +2.50%  '

# Row 12
$ws.Range('B12').Value = 'BinanceUSD'
$ws.Range('C12').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  +0.05%  '

# Row 13
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '22.22'
$ws.Range('E13').Value = '  +3.46%  '

# Row 14
$ws.Range('D14').Value = '6.336'
$ws.Range('E14').Value = '  +4.19%  '

# Row 15
$ws.Range('D15').Value = '7.561'
$ws.Range('E15').Value = '  +5.51%  '

# Row 16
$ws.Range('D16').Value = '1.839.59'
$ws.Range('E16').Value = '  +6.11%  '

# Row 17
$ws.Range('D17').Value = '93.18'
$ws.Range('E17').Value = '  +6.72%  '

# Row 18
$ws.Range('D18').Value = '0.00001086'
$ws.Range('E18').Value = '  +1.87%  '

# Row 19
$ws.Range('D19').Value = '0.06527'
$ws.Range('E19').Value = '  +8.02%  '

# Row 20
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.09%  '

# Row 21
$ws.Range('D21').Value = '17.56'
$ws.Range('E21').Value = '  +4.24%  '

# Row 22
$ws.Range('D22').Value = '6.298'
$ws.Range('E22').Value = '  +3.23%  '

# Row 23
$ws.Range('D23').Value = '28.366.80'

# Row 24
$ws.Range('E24').Value = '  +2.57%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '161.88'
$ws.Range('E25').Value = '  +5.37%  '

# Row 26
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '1.973'
$ws.Range('E26').Value = '  -15.60%  '

# Row 27
$ws.Range('D27').Value = '20.87'
$ws.Range('E27').Value = '  +2.22%  '

# Row 28
$ws.Range('D28').Value = '2.039.89'
$ws.Range('E28').Value = '  +5.46%  '

# Row 29
$ws.Range('D29').Value = '2.303'
$ws.Range('E29').Value = '  -3.28%  '

# Row 30
$ws.Range('D30').Value = '129.26'
$ws.Range('E30').Value = '  +2.81%  '

# Row 31
$ws.Range('D31').Value = '1.226'
$ws.Range('E31').Value = '  +4.12%  '

# Row 32
$ws.Range('D32').Value = '6.045'
$ws.Range('E32').Value = '  +6.34%  '

# Row 33
$ws.Range('D33').Value = '0.09218'
$ws.Range('E33').Value = '  +1.15%  '

# Row 34
$ws.Range('D34').Value = '3.474'
$ws.Range('E34').Value = '  -4.11%  '

# Row 35
$ws.Range('D35').Value = '12.96'
$ws.Range('E35').Value = '  +2.84%  '

# Row 36
$ws.Range('D36').Value = '0.02367'
$ws.Range('E36').Value = '  +3.49%  '

# Row 37
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').Value = '5.241'
$ws.Range('E37').Value = '  +3.42%  '

# Row 38
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').Value = '0.2186'
$ws.Range('E38').Value = '  +2.61%  '

# Row 39
$ws.Range('D39').Value = '0.6612'
$ws.Range('E39').Value = '  +3.67%  '

# Row 40
$ws.Range('D40').Value = '0.06210'
$ws.Range('E40').Value = '  +2.66%  '

# Row 41
$ws.Range('D41').Value = '8.183'
$ws.Range('E41').Value = '  +3.47%  '

# Row 42
$ws.Range('E42').Value = '  +0.34%  '

# Row 43
$ws.Range('D43').Value = '1.440'
$ws.Range('E43').Value = '  +0.95%  '

# Row 44
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  +0.09%  '

# Row 45
$ws.Range('D45').Value = '13.95'
$ws.Range('E45').Value = '  +1.95%  '

# Row 46
$ws.Range('D46').Value = '0.6141'
$ws.Range('E46').Value = '  +5.40%  '

# Row 47
$ws.Range('D47').Value = '3.762'
$ws.Range('E47').Value = '  +1.68%  '

# Row 48
$ws.Range('D48').Value = '126.53'
$ws.Range('E48').Value = '  +0.96%  '

# Row 49
$ws.Range('D49').Value = '2.031'
$ws.Range('E49').Value = '  +4.59%  '

# Row 50
$ws.Range('D50').Value = '1.162'
$ws.Range('E50').Value = '  +5.11%  '

# Row 51
$ws.Range('D51').Value = '0.07022'
$ws.Range('E51').Value = '  +3.00%  '
